$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calice")

# Set new hyperlink text cells
$ws.Range("H19").Value = "CERN_LCDGroup_LucieLinssen_CERN_20140404.pdf"
$ws.Range("H12").Value = "GEM_DHCAL_AndyWhite_UTA_20140326.doc"

# Set G12 date with format
$ws.Range("G12").Value = 41724
$ws.Range("G12").NumberFormat = "d-mmm"

# Move A15:D15 down to A20:D20 (cut/paste, leaving row 15 blank)
$ws.Range("A15:D15").Cut($ws.Range("A20:D20"))

$ws.Range("H16").Select()
